# Update environmental data summary (temp_summary_30_5) and refresh selection.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Nursery)
$ws.Range("B2").Value = 20.197679999999998
$ws.Range("D2").NumberFormat = "0.00"
$ws.Range("D2").Value = 5.9773759999999996
$ws.Range("E2").Value = 50
$ws.Range("F2").Value = 0.68933949999999999
$ws.Range("G2").Value = 11.619619999999999

# Row 3 (CCC)
$ws.Range("D3").NumberFormat = "0.00"
$ws.Range("D3").Value = 3.6226669999999999
$ws.Range("F3").Value = 0.71483540000000001

# Leave the active selection where the author left off editing.
$ws.Range("E11").Select()
